$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-7, columns G through T

$data = @{
    2 = @{ G=8.841467; H=26.524401; I=0.5917001192060068; J=0.5917001192060067; K=3; L=1; M=1.475208; N=4.425624; O=0.8210007041987012; P=0.8210007041987013; Q=13.043002850136; R=117.387025651224; S=0.4857862145425871; T=0.485786214542587 }
    3 = @{ G=8.841467; H=26.524401; I=0.5917001192060068; J=0.5917001192060067; K=3; L=1; M=0.3216333333333333; N=0.9649; O=0.1789992958012987; P=0.1789992958012987; Q=2.843710502766667; R=25.5933945249; S=0.1059139046634197; T=0.1059139046634197 }
    4 = @{ G=4.103438; H=12.310314; I=0.2746155987184545; J=0.2746155987184545; K=3; L=1; M=1.475208; N=4.425624; O=0.8210007041987012; P=0.8210007041987013; Q=6.053424565104; R=54.48082108593599; S=0.2254595999317991; T=0.2254595999317991 }
    5 = @{ G=4.103438; H=12.310314; I=0.2746155987184545; J=0.2746155987184545; K=3; L=1; M=0.3216333333333333; N=0.9649; O=0.1789992958012987; P=0.1789992958012987; Q=1.319802442066667; R=11.8782219786; S=0.04915599878665539; T=0.04915599878665539 }
    6 = @{ G=1.997574666666667; H=5.992724; I=0.1336842820755386; J=0.1336842820755386; K=3; L=1; M=1.475208; N=4.425624; O=0.8210007041987012; P=0.8210007041987013; Q=2.946838128864; R=26.521543159776; S=0.109754889724315; T=0.109754889724315 }
    7 = @{ G=1.997574666666667; H=5.992724; I=0.1336842820755386; J=0.1336842820755386; K=3; L=1; M=0.3216333333333333; N=0.9649; O=0.1789992958012987; P=0.1789992958012987; Q=0.6424865986222222; R=5.7823793876; S=0.02392939235122359; T=0.02392939235122359 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
